# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on row 2 of the
# zh-cn and de-de sheets to reflect the new handback timestamps.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 20:54:35"
$wsZhCn.Range("H2").Value = "2016-03-12 20:54:52"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 20:54:39"
$wsDeDe.Range("H2").Value = "2016-03-12 20:54:58"
